# Update "想去人数" (interested-people count) figures in column F
# across the four worksheets, reflecting refreshed scrape output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 13581
$ws.Range("F6").Value  = 1949
$ws.Range("F9").Value  = 24986
$ws.Range("F10").Value = 554
$ws.Range("F12").Value = 545
$ws.Range("F13").Value = 147
$ws.Range("F14").Value = 394
$ws.Range("F16").Value = 332
$ws.Range("F17").Value = 183
$ws.Range("F18").Value = 156
$ws.Range("F20").Value = 252
$ws.Range("F21").Value = 308
$ws.Range("F22").Value = 33
$ws.Range("F23").Value = 1403
$ws.Range("F24").Value = 110
$ws.Range("F25").Value = 392
$ws.Range("F26").Value = 85

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 4491
$ws.Range("F3").Value  = 209
$ws.Range("F6").Value  = 41
$ws.Range("F10").Value = 404
$ws.Range("F15").Value = 17

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 911
$ws.Range("F3").Value = 4741
$ws.Range("F4").Value = 141

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 911
$ws.Range("F5").Value  = 13581
$ws.Range("F7").Value  = 4741
$ws.Range("F8").Value  = 1949
$ws.Range("F11").Value = 141
$ws.Range("F12").Value = 24986
$ws.Range("F13").Value = 554
$ws.Range("F14").Value = 4491
$ws.Range("F16").Value = 209
$ws.Range("F17").Value = 209
$ws.Range("F18").Value = 545
$ws.Range("F21").Value = 147
$ws.Range("F22").Value = 41
$ws.Range("F26").Value = 404
$ws.Range("F27").Value = 394
$ws.Range("F30").Value = 332
$ws.Range("F31").Value = 183
$ws.Range("F32").Value = 156
$ws.Range("F35").Value = 252
$ws.Range("F38").Value = 308
$ws.Range("F39").Value = 33
$ws.Range("F40").Value = 17
$ws.Range("F41").Value = 1403
$ws.Range("F42").Value = 110
$ws.Range("F44").Value = 392
$ws.Range("F45").Value = 85

$wb.Save()
